# Weekly NYPD CompStat update: new crime data collected.
# Updates report header (volume/number, week-of dates), the weekly crime
# statistics grid (rows 14-31), and the resulting bestFit width of column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a text "placeholder" value (e.g. "0" or "***.*") into a cell
# while preserving that row's normal "label" number format/style (xf 14).
# A leading apostrophe forces Excel to store the value as literal text
# instead of re-parsing "0" as a number; PasteSpecial(formats) afterwards
# re-applies the plain right-aligned text style used elsewhere in the grid
# (anchored on C14, which already carries that exact style).
# ---------------------------------------------------------------------------
function Set-TextPlaceholder($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range("C14").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Report header: "Volume 31   Number  28" -> "...Number  29"
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  29"

# ---------------------------------------------------------------------------
# Report header: week-of date range shifts forward one week
# "Report Covering the Week  7/8/2024  Through  7/14/2024"
#   -> "...7/15/2024  Through  7/21/2024"
# ---------------------------------------------------------------------------
$ws.Range("C9").Value = "Report Covering the Week  7/15/2024  Through  7/21/2024"

# ---------------------------------------------------------------------------
# Crime Complaints grid (rows 14-31): new weekly figures and recomputed
# percentage changes.
# ---------------------------------------------------------------------------

# Row 14 - Murder
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = -66.666666666666

# Row 16 - Robbery
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -87.5
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -58.823529411764
$ws.Range("I16").Value = 70
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = -11.392405063291
$ws.Range("L16").Value = -10.256410256410
$ws.Range("M16").Value = 70.731707317073
$ws.Range("N16").Value = -85.074626865671

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -37.5
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 6.25
$ws.Range("I17").Value = 123
$ws.Range("J17").Value = 75
$ws.Range("K17").Value = 64
$ws.Range("L17").Value = 61.842105263157
$ws.Range("M17").Value = 241.666666666667
$ws.Range("N17").Value = 24.242424242424

# Row 18 - Burglary
Set-TextPlaceholder "C18" "0"
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 9.090909090909
$ws.Range("I18").Value = 131
$ws.Range("J18").Value = 105
$ws.Range("K18").Value = 24.761904761904
$ws.Range("L18").Value = -17.088607594936
$ws.Range("M18").Value = 31
$ws.Range("N18").Value = -71.948608137045

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = -56.25
$ws.Range("F19").Value = 84
$ws.Range("G19").Value = 104
$ws.Range("H19").Value = -19.230769230769
$ws.Range("I19").Value = 613
$ws.Range("J19").Value = 659
$ws.Range("K19").Value = -6.980273141122
$ws.Range("L19").Value = -8.643815201192
$ws.Range("M19").Value = 4.786324786324
$ws.Range("N19").Value = -70.400772573635

# Row 20 - G.L.A.
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 166.666666666667
$ws.Range("J20").Value = 36
$ws.Range("K20").Value = -33.333333333333
$ws.Range("L20").Value = -35.135135135135
$ws.Range("M20").Value = 26.315789473684
$ws.Range("N20").Value = -95.112016293279

# Row 21 - TOTAL
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = -57.692307692307
$ws.Range("F21").Value = 130
$ws.Range("G21").Value = 154
$ws.Range("H21").Value = -15.584415584415
$ws.Range("I21").Value = 973
$ws.Range("J21").Value = 964
$ws.Range("K21").Value = 0.933609958506
$ws.Range("L21").Value = -5.808325266214
$ws.Range("M21").Value = 24.107142857142
$ws.Range("N21").Value = -73.076923076923

# Row 22 - Transit
Set-TextPlaceholder "C22" "0"
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 48
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -21.311475409836

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 86
$ws.Range("D24").Value = 86
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 283
$ws.Range("G24").Value = 311
$ws.Range("H24").Value = -9.003215434083
$ws.Range("I24").Value = 2144
$ws.Range("J24").Value = 2124
$ws.Range("K24").Value = 0.941619585687
$ws.Range("L24").Value = -1.876430205949
$ws.Range("M24").Value = 123.333333333333

# Row 25 - Retail Theft
$ws.Range("C25").Value = 65
$ws.Range("D25").Value = 83
$ws.Range("E25").Value = -21.686746987951
$ws.Range("F25").Value = 258
$ws.Range("G25").Value = 316
$ws.Range("H25").Value = -18.354430379746
$ws.Range("I25").Value = 2088
$ws.Range("J25").Value = 2154
$ws.Range("K25").Value = -3.064066852367
$ws.Range("L25").Value = -5.434782608695

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 28.571428571428
$ws.Range("F26").Value = 34
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 9.677419354838
$ws.Range("I26").Value = 252
$ws.Range("J26").Value = 203
$ws.Range("K26").Value = 24.137931034482
$ws.Range("L26").Value = 29.896907216494
$ws.Range("M26").Value = 80

# Row 28 - Other Sex Crimes
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 62
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = 24
$ws.Range("L28").Value = 6.896551724137

# Row 29 - Shooting Vic.
$ws.Range("M29").Value = 50

# Row 30 - Shooting Inc.
$ws.Range("M30").Value = 50

# Row 31 - Hate Crimes
Set-TextPlaceholder "D31" "0"
Set-TextPlaceholder "E31" "***.*"
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = 0

# ---------------------------------------------------------------------------
# Column H widens (bestFit) to match column E now that it holds a longer
# formatted percentage (e.g. "-58.8"), same width already used by column E.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(5).ColumnWidth
